# Fruta / hortaliza, semanal
# Re-order the weekly price rows (3-13): the D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg) values are
# permuted across rows while all other columns stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: target row -> source row (values that should land in target row come from source row)
$rowMap = @{
    3  = 10
    4  = 3
    5  = 8
    6  = 12
    7  = 4
    8  = 9
    9  = 6
    10 = 7
    11 = 13
    12 = 11
    13 = 5
}

# Snapshot original values for the columns that move, keyed by row number.
$orig = @{}
foreach ($r in 3..13) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $src = $orig[$sourceRow]
    $ws.Cells.Item($targetRow, 4).Value2  = $src.D
    $ws.Cells.Item($targetRow, 10).Value2 = $src.J
    $ws.Cells.Item($targetRow, 11).Value2 = $src.K
    $ws.Cells.Item($targetRow, 12).Value2 = $src.L
    $ws.Cells.Item($targetRow, 13).Value2 = $src.M
    $ws.Cells.Item($targetRow, 16).Value2 = $src.P
}
